# "Generate Report for Archive"
#
# The localization status moved on for the two source files tracked in
# this workbook: they are no longer merely "Ready for handoff" - they are
# now actively "In Translation". That status string is shared across the
# Overview sheet (per-locale status columns) and each per-locale detail
# sheet (zh-cn, de-de), so update it everywhere at once.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# The "Status" column got noticeably narrower once the report was
# regenerated (shorter status text no longer needs as much room). Shrink
# the affected columns to match: columns E/F ("zh-cn"/"de-de" status) on
# the Overview sheet, and column C ("Status") on each locale sheet.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
